$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.278141666666667
$ws.Range("H2").Value = 12.834425
$ws.Range("I2").Value = 0.9663225094340192
$ws.Range("J2").Value = 0.9663225094340191
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.750436666666666
$ws.Range("N2").Value = 14.25131
$ws.Range("O2").Value = 0.07745299862590357
$ws.Range("P2").Value = 0.07745299862590359
$ws.Range("Q2").Value = 20.32304103852778
$ws.Range("R2").Value = 182.90736934675
$ws.Range("S2").Value = 0.07484457599537278
$ws.Range("T2").Value = 0.0748445759953728

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.278141666666667
$ws.Range("H3").Value = 12.834425
$ws.Range("I3").Value = 0.9663225094340192
$ws.Range("J3").Value = 0.9663225094340191
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 34.05277366666667
$ws.Range("N3").Value = 102.158321
$ws.Range("O3").Value = 0.5552098927072401
$ws.Range("P3").Value = 0.5552098927072401
$ws.Range("Q3").Value = 145.6825898889361
$ws.Range("R3").Value = 1311.143309000425
$ws.Range("S3").Value = 0.5365118167834528
$ws.Range("T3").Value = 0.5365118167834527

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.278141666666667
$ws.Range("H4").Value = 12.834425
$ws.Range("I4").Value = 0.9663225094340192
$ws.Range("J4").Value = 0.9663225094340191
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 22.52994333333334
$ws.Range("N4").Value = 67.58983
$ws.Range("O4").Value = 0.3673371086668564
$ws.Range("P4").Value = 0.3673371086668564
$ws.Range("Q4").Value = 96.38628932197223
$ws.Range("R4").Value = 867.47660389775
$ws.Range("S4").Value = 0.3549661166551937
$ws.Range("T4").Value = 0.3549661166551937

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.1490983333333333
$ws.Range("H5").Value = 0.447295
$ws.Range("I5").Value = 0.03367749056598092
$ws.Range("J5").Value = 0.03367749056598091
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.750436666666666
$ws.Range("N5").Value = 14.25131
$ws.Range("O5").Value = 0.07745299862590357
$ws.Range("P5").Value = 0.07745299862590359
$ws.Range("Q5").Value = 0.7082821896055556
$ws.Range("R5").Value = 6.37453970645
$ws.Range("S5").Value = 0.002608422630530801
$ws.Range("T5").Value = 0.002608422630530801

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.1490983333333333
$ws.Range("H6").Value = 0.447295
$ws.Range("I6").Value = 0.03367749056598092
$ws.Range("J6").Value = 0.03367749056598091
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 34.05277366666667
$ws.Range("N6").Value = 102.158321
$ws.Range("O6").Value = 0.5552098927072401
$ws.Range("P6").Value = 0.5552098927072401
$ws.Range("Q6").Value = 5.077211799077222
$ws.Range("R6").Value = 45.694906191695
$ws.Range("S6").Value = 0.01869807592378736
$ws.Range("T6").Value = 0.01869807592378735

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.1490983333333333
$ws.Range("H7").Value = 0.447295
$ws.Range("I7").Value = 0.03367749056598092
$ws.Range("J7").Value = 0.03367749056598091
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 22.52994333333334
$ws.Range("N7").Value = 67.58983
$ws.Range("O7").Value = 0.3673371086668564
$ws.Range("P7").Value = 0.3673371086668564
$ws.Range("Q7").Value = 3.359177001094445
$ws.Range("R7").Value = 30.23259300985
$ws.Range("S7").Value = 0.01237099201166276
$ws.Range("T7").Value = 0.01237099201166276
